# Update countries & provincias Spain
# Applies the data refresh captured in the upstream OOXML diff:
#   - bump the "Datos actualizados ..." timestamp string
#   - swap the row labels for Nepal/Uzbekistan (Uzbekistan now outranks
#     Nepal in total cases, so their shared-string-backed labels trade
#     places while the numeric rows keep their own new figures)
#   - swap the row labels for Islas Malvinas/Groenlandia likewise
#   - refresh the numeric counters (Total, Nuevos, Activos, Recuperados,
#     Criticos, Muertes) for the affected country rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 22 de Julio de 2020 a las 07:51"

# --- Row 15: Pakistan ---------------------------------------------------
$ws.Range("B15").Value = 267428
$ws.Range("C15").Value = 1332
$ws.Range("D15").Value = 210468
$ws.Range("E15").Value = 51283
$ws.Range("G15").Value = 38
$ws.Range("H15").Value = 5677

# --- Row 57: Kirguistan --------------------------------------------------
$ws.Range("B57").Value = 28980
$ws.Range("C57").Value = 729
$ws.Range("D57").Value = 15536
$ws.Range("E57").Value = 12333
$ws.Range("G57").Value = 32
$ws.Range("H57").Value = 1111

# --- Rows 65/66: Nepal & Uzbekistan swap ranking -------------------------
# Row 65 becomes Uzbekistan with its refreshed figures.
$ws.Range("A65").Value = "Uzbekistan"
$ws.Range("B65").Value = 18171
$ws.Range("C65").Value = 290
$ws.Range("D65").Value = 9521
$ws.Range("E65").Value = 8554
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 96

# Row 66 becomes Nepal, carrying the figures Nepal previously had at row 65.
$ws.Range("A66").Value = "Nepal"
$ws.Range("B66").Value = 17994
$ws.Range("D66").Value = 12477
$ws.Range("E66").Value = 5477
$ws.Range("H66").Value = 40

# --- Row 90: Haiti ---------------------------------------------------
$ws.Range("B90").Value = 7146
$ws.Range("C90").Value = 46
$ws.Range("E90").Value = 2897
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 154

# --- Row 106: Tailandia --------------------------------------------------
$ws.Range("B106").Value = 3261
$ws.Range("C106").Value = 6
$ws.Range("E106").Value = 98

# --- Row 196: Belice -------------------------------------------------
$ws.Range("B196").Value = 43
$ws.Range("C196").Value = 1
$ws.Range("D196").Value = 23

# --- Row 200: Papua Nueva Guinea -----------------------------------------
$ws.Range("E200").Value = 19
$ws.Range("H200").Value = 0

# --- Rows 210/211: Islas Malvinas & Groenlandia swap ranking -------------
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"
